# daily auto push: 2025-09-30 01:58 UTC
# Append the new daily-tracking row (row 39) to Sheet1:
#   A39 = "2025/09/30" (text, same style as the date column above)
#   B39 = "火"          (text, day-of-week)
#   C39 = 9             (number, hour)
#   D39 = 16            (number, ranking)
# The sheet's dimension (A1:D38 -> A1:D39) updates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Cells.Item(39, 1)
# Force the date-like string to be stored as literal text (matching the
# existing rows), instead of letting Excel auto-convert it to a date
# serial number. Clear the temporary number-format afterwards so the new
# cell ends up with no explicit style, just like its neighbours.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025/09/30"
$dateCell.ClearFormats()

$ws.Cells.Item(39, 2).Value = "火"
$ws.Cells.Item(39, 3).Value = 9
$ws.Cells.Item(39, 4).Value = 16
